# Summary of score frequencies E1 and E2
#
# This script reproduces the authoring diff:
#  - "Sheet1" becomes "Score" and is populated with new score-frequency data,
#    formulas and a table.
#  - "Sheet2" becomes "Gender by Profession" (chart series formulas updated
#    to point at the new sheet name).
#  - Workbook view: Score tab becomes active, first visible tab shifts.
#  - The unused external link (gender_profession.csv) is broken/removed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update chart series formulas on "Sheet2" BEFORE renaming it, so we
#    can match on the old sheet name text inside the SERIES() formula.
# ---------------------------------------------------------------------
$sheet2 = $wb.Worksheets.Item("Sheet2")
for ($i = 1; $i -le $sheet2.ChartObjects().Count; $i++) {
    $co = $sheet2.ChartObjects($i)
    $chart = $co.Chart
    for ($j = 1; $j -le $chart.SeriesCollection().Count; $j++) {
        $ser = $chart.SeriesCollection($j)
        $ser.Formula = $ser.Formula.Replace("Sheet2!", "'Gender by Profession'!")
    }
}

# ---------------------------------------------------------------------
# 2. Rename the sheets.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Sheet1").Name = "Score"
$sheet2.Name = "Gender by Profession"

# ---------------------------------------------------------------------
# 3. Break the now-unused external link to gender_profession.csv.
# ---------------------------------------------------------------------
$links = $wb.LinkSources(1)
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}

# ---------------------------------------------------------------------
# 4. Populate the "Score" sheet with the new content.
# ---------------------------------------------------------------------
$score = $wb.Worksheets.Item("Score")
$score.Range("A1").ClearContents()

$score.Range("A2").Value = "Column1"
$score.Range("B2").Value = "low"
$score.Range("C2").Value = "medium"
$score.Range("D2").Value = "high"
$score.Range("E2").Value = "total"

$score.Range("A3").Value = "E1"
$score.Range("B3").Value = 538
$score.Range("C3").Value = 134
$score.Range("D3").Value = 105
$score.Range("E3").Formula = "=SUM(B3:D3)"

$score.Range("A4").Value = "E2"
$score.Range("B4").Value = 146
$score.Range("C4").Value = 157
$score.Range("D4").Value = 194
$score.Range("E4").Formula = "=SUM(B4:D4)"

$score.Range("A5").Value = "total"
$score.Range("B5").Formula = "=SUM(B3:B4)"
$score.Range("C5").Formula = "=SUM(C3:C4)"
$score.Range("D5").Formula = "=SUM(D3:D4)"

$score.Range("A7").Value = "Hobbyist"
$score.Range("B7").Value = "low"
$score.Range("C7").Value = "medium"
$score.Range("D7").Value = "high"

$score.Range("A8").Value = "E1"
$score.Range("B8").Formula = "=B3/B5"
$score.Range("C8").Formula = "=C3/C5"
$score.Range("D8").Formula = "=D3/D5"

$score.Range("A9").Value = "E2"
$score.Range("B9").Formula = "=B4/B5"
$score.Range("C9").Formula = "=C4/C5"
$score.Range("D9").Formula = "=D4/D5"

$score.Range("A11").Value = "This shows that E2 qualification test spread participants more evenly across the three qualification levels."
$score.Range("A12").Value = "The Chisquare test confirmed that that experiment assignment and the levels of qualification were independent."

$score.Range("A14").Value = "We donfirmed that these proportion are distinct by running a chisquare test to evaluate the independence"
$score.Range("A15").Value = "between independent from the qualification score levels (chisquare 201,14, df=2, p-value<0.0001)"

# Column widths
$score.Columns.Item(1).ColumnWidth = 9.86328125
$score.Columns.Item(3).ColumnWidth = 9.3984375

# Table over the frequency block.
$score.ListObjects.Add(1, $score.Range("A2:E5"), 0, 1).Name = "Table11"

# ---------------------------------------------------------------------
# 5. Selection / view bookkeeping matching the diff.
# ---------------------------------------------------------------------
$score.Range("A17").Select()

$gbp = $wb.Worksheets.Item("Gender by Profession")
$gbp.Range("G13").Select()

$score.Activate()
$wb.Windows.Item(1).DisplayedTabs = 1
